$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Septiembre de 2020 a las 02:23"

# 2) Swap country-name labels for the four re-sorted pairs (values unchanged unless noted below)
$ws.Cells.Item(61,1).Value = "Chequia"
$ws.Cells.Item(62,1).Value = "Suiza"
$ws.Cells.Item(126,1).Value = "Surinam"
$ws.Cells.Item(127,1).Value = "Ruanda"
$ws.Cells.Item(204,1).Value = "Timor Oriental"
$ws.Cells.Item(205,1).Value = "Santa Lucia"
$ws.Cells.Item(214,1).Value = "Islas Malvinas"
$ws.Cells.Item(215,1).Value = "Montserrat"

# 3) Update changed numeric statistics per row
$ws.Cells.Item(4,2).Value = 7001273
$ws.Cells.Item(4,3).Value = 32894
$ws.Cells.Item(4,4).Value = 4248346
$ws.Cells.Item(4,5).Value = 2548809
$ws.Cells.Item(4,7).Value = 294
$ws.Cells.Item(4,8).Value = 204118

$ws.Cells.Item(25,4).Value = 244000
$ws.Cells.Item(25,5).Value = 20007

$ws.Cells.Item(29,2).Value = 143651
$ws.Cells.Item(29,3).Value = 877
$ws.Cells.Item(29,4).Value = 124691
$ws.Cells.Item(29,5).Value = 9743

$ws.Cells.Item(61,2).Value = 49290
$ws.Cells.Item(61,3).Value = 984
$ws.Cells.Item(61,4).Value = 24755
$ws.Cells.Item(61,5).Value = 24032
$ws.Cells.Item(61,7).Value = 4
$ws.Cells.Item(61,8).Value = 503

$ws.Cells.Item(62,2).Value = 49283
$ws.Cells.Item(62,3).Value = 0
$ws.Cells.Item(62,4).Value = 40500
$ws.Cells.Item(62,5).Value = 6738
$ws.Cells.Item(62,7).Value = 0
$ws.Cells.Item(62,8).Value = 2045

$ws.Cells.Item(72,2).Value = 33520
$ws.Cells.Item(72,3).Value = 505
$ws.Cells.Item(72,4).Value = 18117
$ws.Cells.Item(72,5).Value = 14744
$ws.Cells.Item(72,7).Value = 23
$ws.Cells.Item(72,8).Value = 659

$ws.Cells.Item(92,2).Value = 13555
$ws.Cells.Item(92,3).Value = 20
$ws.Cells.Item(92,4).Value = 6760
$ws.Cells.Item(92,5).Value = 5959

$ws.Cells.Item(93,2).Value = 12897
$ws.Cells.Item(93,3).Value = 39
$ws.Cells.Item(93,5).Value = 2259

$ws.Cells.Item(105,2).Value = 8619
$ws.Cells.Item(105,3).Value = 4
$ws.Cells.Item(105,5).Value = 2035

$ws.Cells.Item(106,2).Value = 8612
$ws.Cells.Item(106,3).Value = 714
$ws.Cells.Item(106,4).Value = 5268
$ws.Cells.Item(106,5).Value = 3208
$ws.Cells.Item(106,7).Value = 2
$ws.Cells.Item(106,8).Value = 136

$ws.Cells.Item(107,2).Value = 7907
$ws.Cells.Item(107,3).Value = 103
$ws.Cells.Item(107,5).Value = 1080

$ws.Cells.Item(108,2).Value = 7683
$ws.Cells.Item(108,3).Value = 11
$ws.Cells.Item(108,4).Value = 5924
$ws.Cells.Item(108,5).Value = 1534

$ws.Cells.Item(115,4).Value = 5335
$ws.Cells.Item(115,5).Value = 7

$ws.Cells.Item(116,2).Value = 5269
$ws.Cells.Item(116,3).Value = 24
$ws.Cells.Item(116,4).Value = 4624
$ws.Cells.Item(116,5).Value = 541

$ws.Cells.Item(117,2).Value = 5257
$ws.Cells.Item(117,3).Value = 71
$ws.Cells.Item(117,4).Value = 4599
$ws.Cells.Item(117,5).Value = 607
$ws.Cells.Item(117,7).Value = 1
$ws.Cells.Item(117,8).Value = 51

$ws.Cells.Item(126,2).Value = 4723
$ws.Cells.Item(126,3).Value = 14
$ws.Cells.Item(126,4).Value = 4488
$ws.Cells.Item(126,5).Value = 138
$ws.Cells.Item(126,8).Value = 97

$ws.Cells.Item(127,2).Value = 4711
$ws.Cells.Item(127,3).Value = 22
$ws.Cells.Item(127,4).Value = 2961
$ws.Cells.Item(127,5).Value = 1724
$ws.Cells.Item(127,8).Value = 26

$ws.Cells.Item(129,2).Value = 3991
$ws.Cells.Item(129,3).Value = 90
$ws.Cells.Item(129,5).Value = 2394
$ws.Cells.Item(129,7).Value = 5
$ws.Cells.Item(129,8).Value = 152

$ws.Cells.Item(151,2).Value = 2269
$ws.Cells.Item(151,3).Value = 101
$ws.Cells.Item(151,4).Value = 1339
$ws.Cells.Item(151,5).Value = 866

$ws.Cells.Item(155,2).Value = 1846
$ws.Cells.Item(155,3).Value = 30
$ws.Cells.Item(155,4).Value = 1187
$ws.Cells.Item(155,5).Value = 603

$ws.Cells.Item(159,2).Value = 1600
$ws.Cells.Item(159,3).Value = 10
$ws.Cells.Item(159,5).Value = 296

$ws.Cells.Item(214,4).Value = 13
$ws.Cells.Item(214,8).Value = 0

$ws.Cells.Item(215,4).Value = 12
$ws.Cells.Item(215,8).Value = 1

